$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

# New row 4 mirrors the layout of row 3 (FunNm / key-read condition / other order condition)
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(4, 2).Value = "EntryDate >= ,AND EntryDate <= ,AND Factor = , AND CustNo = "
$ws.Cells.Item(4, 1).Value = "findEntryDateRangeFactorCustNoFirst"
$ws.Cells.Item(4, 3).Value = "EntryDate"

# Update the remembered selection on DBD (no longer the active tab)
$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBD.Activate() | Out-Null
$wsDBD.Range("C9").Select() | Out-Null

# DBS becomes the active sheet with C5 selected
$ws.Activate() | Out-Null
$ws.Range("C5").Select() | Out-Null
